# Estadisticos Matutinos 15 Oct
$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 7
$ws1.Range("F2").Value = 24
$ws1.Range("G2").Value = 77.42
$ws1.Range("H2").Value = 8

$ws1.Range("D3").Value = 10
$ws1.Range("F3").Value = 11
$ws1.Range("G3").Value = 52.38
$ws1.Range("H3").Value = 8

$ws1.Range("D4").Value = 7
$ws1.Range("F4").Value = 28
$ws1.Range("G4").Value = 80

$ws1.Range("D5").Value = 8
$ws1.Range("F5").Value = 13
$ws1.Range("G5").Value = 61.9
$ws1.Range("H5").Value = 7.8

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 24
$ws2.Range("E3").Value = 11
$ws2.Range("E4").Value = 28
$ws2.Range("E5").Value = 13

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 7
$ws3.Range("F2").Value = 24
$ws3.Range("G2").Value = 77.42
$ws3.Range("H2").Value = 8

$ws3.Range("D3").Value = 10
$ws3.Range("F3").Value = 11
$ws3.Range("G3").Value = 52.38
$ws3.Range("H3").Value = 8

$ws3.Range("D4").Value = 7
$ws3.Range("F4").Value = 28
$ws3.Range("G4").Value = 80

$ws3.Range("D5").Value = 8
$ws3.Range("F5").Value = 13
$ws3.Range("G5").Value = 61.9
$ws3.Range("H5").Value = 7.8

# --- Sheet "Rescatables" ---
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Cells.Item(2, 1).Value = 20330051920396
$ws4.Cells.Item(2, 2).Value = "MUÑOZ"
$ws4.Cells.Item(2, 3).Value = "RODRIGUEZ"
$ws4.Cells.Item(2, 4).Value = "VICTOR HUGO"
$ws4.Cells.Item(2, 5).Value = "GEOMETRÍA ANALÍTICA"
$ws4.Cells.Item(2, 6).Value = "3ASV"
$ws4.Cells.Item(2, 7).Value = 6

$ws4.Cells.Item(3, 1).Value = 20330051920119
$ws4.Cells.Item(3, 2).Value = "CORONA"
$ws4.Cells.Item(3, 3).Value = "HERNANDEZ"
$ws4.Cells.Item(3, 4).Value = "MARIA FERNANDA"
$ws4.Cells.Item(3, 5).Value = "GEOMETRÍA ANALÍTICA"
$ws4.Cells.Item(3, 6).Value = "3ARHV"
$ws4.Cells.Item(3, 7).Value = 6

$ws4.Cells.Item(4, 1).Value = 20330051920128
$ws4.Cells.Item(4, 2).Value = "HERNANDEZ"
$ws4.Cells.Item(4, 3).Value = "SUAREZ"
$ws4.Cells.Item(4, 4).Value = "KIMBERLY ALONDRA"
$ws4.Cells.Item(4, 5).Value = "GEOMETRÍA ANALÍTICA"
$ws4.Cells.Item(4, 6).Value = "3ARHV"
$ws4.Cells.Item(4, 7).Value = 6
